$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Phi_phi0"
$ws.Range("N1").Value = "Phi_theta0"
$ws.Range("B2").Value = 0.0005865707428022582
$ws.Range("C2").Value = 0.00000003101067652844651
$ws.Range("D2").Value = 0.000000000476923237286854
$ws.Range("E2").Value = 0.000000001217627311888573
$ws.Range("F2").Value = -0.000000002149355397241235
$ws.Range("G2").Value = -0.000000002943145147442055
$ws.Range("H2").Value = -0.00002207081446607662
$ws.Range("I2").Value = -0.00001364189471955363
$ws.Range("J2").Value = 0.0000223736857949539
$ws.Range("K2").Value = -0.00000993891657917283
$ws.Range("L2").Value = 0.000006683802486366034
$ws.Range("M2").Value = 0.00008610866161386833
$ws.Range("N2").Value = 0.0001217067026377315
$ws.Range("O2").Value = 0.0002358803922126165
$ws.Range("C3").Value = 0.0000000000612233567191681
$ws.Range("D3").Value = 0.00000000001623206650674353
$ws.Range("E3").Value = -0.000000000006363829484292395
$ws.Range("F3").Value = 0.000000000004270705028095265
$ws.Range("G3").Value = 0.000000000002809788026041792
$ws.Range("H3").Value = -0.00000001932949147894004
$ws.Range("I3").Value = 0.000000003298677412050234
$ws.Range("J3").Value = -0.000000001322792342254664
$ws.Range("K3").Value = -0.000000006985411301524843
$ws.Range("L3").Value = 0.000000003283835066933199
$ws.Range("M3").Value = -0.00000009630816343488513
$ws.Range("N3").Value = 0.0000000350994249248629
$ws.Range("O3").Value = -0.0000001450560043107301
$ws.Range("D4").Value = 0.00000000004846513084217753
$ws.Range("E4").Value = -0.0000000000003060066258302464
$ws.Range("F4").Value = 0.0000000000002192492027873383
$ws.Range("G4").Value = -0.00000000003902789762772377
$ws.Range("H4").Value = 0.000000002996140051948911
$ws.Range("I4").Value = 0.000000003633962267949914
$ws.Range("J4").Value = -0.0000000009374353840181265
$ws.Range("K4").Value = 0.00000001078494016713809
$ws.Range("L4").Value = 0.00000001708181586882974
$ws.Range("M4").Value = -0.0000002653782791289729
$ws.Range("N4").Value = 0.0000002285532895614608
$ws.Range("O4").Value = -0.00000004787361765700931
$ws.Range("E5").Value = 0.000000000003362010712423705
$ws.Range("F5").Value = -0.000000000002380107134266269
$ws.Range("G5").Value = -0.000000000003124884629242691
$ws.Range("H5").Value = -0.000000006265488494255017
$ws.Range("I5").Value = -0.0000000003210061486853674
$ws.Range("J5").Value = 0.0000000004355921438090672
$ws.Range("K5").Value = 0.0000000007776594098431285
$ws.Range("L5").Value = -0.000000003077394067943565
$ws.Range("M5").Value = -0.00000002714223289424777
$ws.Range("N5").Value = -0.0000000211591505308449
$ws.Range("O5").Value = 0.00000002846239149653389
$ws.Range("F6").Value = 0.000000000002614607825925151
$ws.Range("G6").Value = 0.000000000001623029747288337
$ws.Range("H6").Value = 0.00000002807436162137531
$ws.Range("I6").Value = -0.0000000002878306849915708
$ws.Range("J6").Value = -0.0000000009690951294949218
$ws.Range("K6").Value = 0.000000000256476600227827
$ws.Range("L6").Value = -0.000000000891299310853885
$ws.Range("M6").Value = 0.0000000473089343297891
$ws.Range("N6").Value = 0.00000004191687930943955
$ws.Range("O6").Value = 0.000000002276717957602607
$ws.Range("G7").Value = 0.00000000004166290235743476
$ws.Range("H7").Value = 0.000000001453144874452881
$ws.Range("I7").Value = -0.000000002258239862054472
$ws.Range("J7").Value = 0.0000000004423353537295178
$ws.Range("K7").Value = -0.00000001841641977164669
$ws.Range("L7").Value = -0.00000003354091178562904
$ws.Range("M7").Value = 0.000000294562800927592
$ws.Range("N7").Value = -0.0000001081899147480721
$ws.Range("O7").Value = -0.00000008413375929841176
$ws.Range("H8").Value = 0.002422240117188004
$ws.Range("I8").Value = -0.000022147135736164
$ws.Range("J8").Value = -0.00001783004619602197
$ws.Range("K8").Value = 0.00009729648590519111
$ws.Range("L8").Value = 0.0001440123414253604
$ws.Range("M8").Value = 0.001487490293657892
$ws.Range("N8").Value = 0.00125524280802329
$ws.Range("O8").Value = -0.0003142943039060382
$ws.Range("I9").Value = 0.002232303694120674
$ws.Range("J9").Value = 0.00002165489783586008
$ws.Range("K9").Value = -0.001030141029337531
$ws.Range("L9").Value = 0.0008220734456962875
$ws.Range("M9").Value = 0.000004293781999946261
$ws.Range("N9").Value = 0.0001879391457534538
$ws.Range("O9").Value = 0.0001178874965792063
$ws.Range("J10").Value = 0.0003533587771233257
$ws.Range("K10").Value = -0.00001615476956471619
$ws.Range("L10").Value = 0.0002760442350934531
$ws.Range("M10").Value = -0.0000899426471558907
$ws.Range("N10").Value = -0.000009827263126949856
$ws.Range("O10").Value = -0.00000586232975499785
$ws.Range("K11").Value = 0.02656093706752399
$ws.Range("L11").Value = 0.01207856161886289
$ws.Range("M11").Value = -0.0129523303797163
$ws.Range("N11").Value = 0.1173198877868248
$ws.Range("O11").Value = -0.0001247682940825344
$ws.Range("L12").Value = 0.1690812681821677
$ws.Range("M12").Value = -0.1200937770874017
$ws.Range("N12").Value = 0.1024552451204341
$ws.Range("O12").Value = 0.0004649557183792104
$ws.Range("A13").Value = "Phi_phi0"
$ws.Range("M13").Value = 0.1986720829524958
$ws.Range("N13").Value = -0.06162235256176028
$ws.Range("O13").Value = -0.005354331767963187
$ws.Range("A14").Value = "Phi_theta0"
$ws.Range("N14").Value = 1.056799840188618
$ws.Range("O14").Value = -0.005290252050545846
$ws.Range("O15").Value = 0.0479951001769764
